$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.984.36"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.829.78"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'706.90"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "'171.44"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "3.828.89"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "4.474.85"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "3.802.09"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").Value = "70.998.29"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'7.22"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'17.37"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'495.21"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "'10.64"
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("D23").Value = "'0.734"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "'85.43"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").Value = "3.982.69"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D31").Value = "'3.10"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'7.41"
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("D34").Value = "'29.33"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.799.44"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'9.16"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "'5.97"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").Value = "'3.31"
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'163.67"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'428.82"
$ws.Range("E48").Value = "  +3.67%  "
$ws.Range("D49").Value = "'48.92"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'8.74"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'1.36"
$ws.Range("E51").Value = "  -1.58%  "
